$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings as PS variables (indices 24-48) ---
$s24 = "SCRIPT/T01P01A/um1308.ssb"
$s25 = " I will lend my support to\ncapturing that thieving [CS:N]Grovyle[CR]."
$s26 = " Let us work together...[K]and catch\nthat [CS:N]Grovyle[CR]."
$s27 = " Ÿ ðïíïãô âàí â ðïéíëå üóïãï\nðïöéóéóåìÿ [CS:N]Ãñïâàêìà[CR]."
$s28 = " Ðïóñôäéíòÿ òïïáþà...[K] É ðïêíàåí\n[CS:N]Ãñïâàêìà[CR]."
$s29 = " Я помогу вам в поимке этого\nпохитителя [CS:N]Гровайла[CR]."
$s30 = " Потрудимся сообща...[K] И поймаем\n[CS:N]Гровайла[CR]."
$s31 = " The [CS:P]Northern Desert[CR] is vast and\ndeep. It is also frequently scoured by savage\nsandstorms."
$s32 = " Please do take care!"
$s33 = " [CS:P]Северная Пустыня[CR] обширна,\nглубока и к тому же, в ней часто возникают\nсильные песчаные бури."
$s34 = " [CS:P]Òåâåñîàÿ Ðôòóúîÿ[CR] ïášéñîà,\nãìôáïëà é ë óïíô çå, â îåê œàòóï âïèîéëàýó\nòéìûîúå ðåòœàîúå áôñé."
$s35 = " Пожалуйста, берегите себя!"
$s36 = " Ðïçàìôêòóà, áåñåãéóå òåáÿ!"
$s37 = "SCRIPT/G01P03A/um1307.ssb"
$s38 = "SCRIPT/G01P04A/um1303.ssb"
$s39 = " I'm sorry. I'm certain that we'll\ndevise our next plan soon."
$s40 = " Until then, please conduct\nsearches on your own."
$s41 = " That's all I can ask of you now."
$s42 = " Мне очень жаль. Я уверен, что\nскоро мы придумаем новый план."
$s43 = " До сей поры, вам придётся\nискать самим."
$s44 = " Это всё, о чем я могу вас\nпросить."
$s45 = " Íîå ïœåîû çàìû. Ÿ ôâåñåî, œóï\nòëïñï íú ðñéäôíàåí îïâúê ðìàî."
$s46 = " Äï òåê ðïñú, âàí ðñéäæóòÿ\néòëàóû òàíéí."
$s47 = " Üóï âòæ, ï œæí ÿ íïãô âàò\nðñïòéóû."
$s48 = "SCRIPT/G01P04A/um1402.ssb"
# --- Row 7: turn into the last (bordered) row of the first (pre-existing) group ---
$ws.Range("A7").WrapText = $true
$b = $ws.Range("A7").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

$ws.Range("B7").WrapText = $true
$b = $ws.Range("B7").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

$ws.Range("C7").WrapText = $true
$ws.Range("C7").Font.Size = 8
$b = $ws.Range("C7").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

$ws.Range("D7").WrapText = $true
$ws.Range("D7").Font.Size = 8
$b = $ws.Range("D7").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

$ws.Range("E7").WrapText = $true
$ws.Range("E7").Font.Size = 8
$b = $ws.Range("E7").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

# --- Group 1 (rows 8-9): "SCRIPT/T01P01A/um1308.ssb" ---
# Cell values are entered in the exact order needed to reproduce the
# original shared-string table ordering.
$ws.Range("A8").Value = $s24
$ws.Range("C8").Value = $s25
$ws.Range("C9").Value = $s26
$ws.Range("E8").Value = $s27
$ws.Range("E9").Value = $s28
$ws.Range("D8").Value = $s29
$ws.Range("D9").Value = $s30
$ws.Range("B8").Value = 66
$ws.Range("B9").Value = 69

$ws.Rows.Item(8).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 21.6

$ws.Range("A9").WrapText = $true
$b = $ws.Range("A9").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

$ws.Range("B9").WrapText = $true
$b = $ws.Range("B9").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

$ws.Range("C9").WrapText = $true
$ws.Range("C9").Font.Size = 8
$b = $ws.Range("C9").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

$ws.Range("D9").WrapText = $true
$ws.Range("D9").Font.Size = 8
$b = $ws.Range("D9").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

$ws.Range("E9").WrapText = $true
$ws.Range("E9").Font.Size = 8
$b = $ws.Range("E9").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

# --- Group 2 (rows 10-11): "SCRIPT/G01P03A/um1307.ssb" ---
$ws.Range("C10").Value = $s31
$ws.Range("C11").Value = $s32
$ws.Range("D10").Value = $s33
$ws.Range("E10").Value = $s34
$ws.Range("D11").Value = $s35
$ws.Range("E11").Value = $s36
$ws.Range("A10").Value = $s37
$ws.Range("B10").Value = 44
$ws.Range("B11").Value = 47

$ws.Rows.Item(10).RowHeight = 43.2

$ws.Range("D10").WrapText = $true
$ws.Range("D10").Font.Size = 8
$ws.Range("E10").WrapText = $true
$ws.Range("E10").Font.Size = 8
$ws.Range("F10").WrapText = $true
$ws.Range("F11").WrapText = $true

$ws.Range("A11").WrapText = $true
$b = $ws.Range("A11").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

$ws.Range("B11").WrapText = $true
$b = $ws.Range("B11").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

$ws.Range("C11").WrapText = $true
$ws.Range("C11").Font.Size = 8
$b = $ws.Range("C11").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

$ws.Range("D11").WrapText = $true
$ws.Range("D11").Font.Size = 8
$b = $ws.Range("D11").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

$ws.Range("E11").WrapText = $true
$ws.Range("E11").Font.Size = 8
$b = $ws.Range("E11").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2

# --- Group 3 (rows 12-14): "SCRIPT/G01P04A/um1303.ssb" + "SCRIPT/G01P04A/um1402.ssb" ---
$ws.Range("A12").Value = $s38
$ws.Range("C12").Value = $s39
$ws.Range("C13").Value = $s40
$ws.Range("C14").Value = $s41
$ws.Range("D12").Value = $s42
$ws.Range("D13").Value = $s43
$ws.Range("D14").Value = $s44
$ws.Range("E12").Value = $s45
$ws.Range("E13").Value = $s46
$ws.Range("E14").Value = $s47
$ws.Range("A13").Value = $s48
$ws.Range("B12").Value = 18
$ws.Range("B13").Value = 21
$ws.Range("B14").Value = 24

$ws.Rows.Item(12).RowHeight = 43.2
$ws.Rows.Item(13).RowHeight = 43.2

# --- View state: selection / scroll position ---
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("D15").Select()
